$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.883.20'
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = '2.482.85'
$ws.Range("E3").Value = '  -1.38%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("E5").Value = '  -1.49%  '

$ws.Range("D6").Formula = "=""104.23"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -4.71%  '

$ws.Range("E7").Value = '  -2.35%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").Formula = "=""0.534"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -2.73%  '

$ws.Range("D10").Formula = "=""38.63"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -3.93%  '

$ws.Range("D11").Formula = "=""20.29"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").Formula = "=""0.0797"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -2.77%  '

$ws.Range("E13").Value = '  +0.51%  '

$ws.Range("E14").Value = '  -2.96%  '

$ws.Range("D15").Value = '2.873.64'
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").Value = '2.484.27'
$ws.Range("E16").Value = '  -1.70%  '

$ws.Range("E17").Value = '  -3.28%  '

$ws.Range("D18").Value = '47.788.17'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("E19").Value = '  +8.12%  '

$ws.Range("E20").Value = '  -4.78%  '

$ws.Range("D21").Formula = "=""6.51"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = '0.0₃0925'
$ws.Range("E22").Value = '  -2.23%  '

$ws.Range("D23").Formula = "=""278.30"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +1.80%  '

$ws.Range("D24").Formula = "=""70.60"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -1.91%  '

$ws.Range("E25").Value = '  -3.82%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").Formula = "=""25.53"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -1.62%  '

$ws.Range("D28").Formula = "=""2.24"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -5.77%  '

$ws.Range("E29").Value = '  -5.33%  '

$ws.Range("E30").Value = '  -4.68%  '

$ws.Range("E31").Value = '  -3.47%  '

$ws.Range("D32").Formula = "=""48.99"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -1.27%  '

$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").Formula = "=""18.96"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -3.61%  '

$ws.Range("D35").Formula = "=""5.22"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("D36").Formula = "=""0.0767"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -2.19%  '

$ws.Range("E37").Value = '  -2.66%  '

$ws.Range("E38").Value = '  -4.38%  '

$ws.Range("E39").Value = '  -4.84%  '

$ws.Range("E40").Value = '  -1.23%  '

$ws.Range("D41").Formula = "=""2.20"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -1.75%  '

$ws.Range("D42").Formula = "=""119.71"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -1.57%  '

$ws.Range("D43").Formula = "=""21.41"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -2.38%  '

$ws.Range("D44").Formula = "=""0.0297"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)

$ws.Range("D45").Value = '1.985.06'
$ws.Range("E45").Value = '  -2.20%  '

$ws.Range("D46").Formula = "=""3.09"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -1.42%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Formula = "=""2.09"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Formula = "=""1.90"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("D49").Formula = "=""8.88"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -2.26%  '

$ws.Range("E50").Value = '  -2.43%  '

$ws.Range("D51").Formula = "=""78.53"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -1.86%  '

$excel.CutCopyMode = $false